$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ризики")
$ws.Activate()

$ws.Range("J2").Value = "Артем Ліфтієв"
$ws.Range("J3").Value = "Артем Ліфтієв"
$ws.Range("J4").Value = "Павло Луговий"
$ws.Range("J5").Value = "Таїсія Деркач"
$ws.Range("J6").Value = "Леонід Головненко"
$ws.Range("J7").Value = "Іван Кражан"
$ws.Range("J8").Value = "Артем Ліфтієв"
$ws.Range("J9").Value = "Анна Кравченко"

$ws.Range("B7").Value = "Технічні проблеми з ігровим двигуном"
$ws.Range("B7").Font.Name = "Arial"

$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("J5").Select()
